$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 0.5903170108795166
$ws.Cells.Item(2, 5).Value = 4534.651566234569
$ws.Cells.Item(2, 6).Value = 0.1022928192822754
$ws.Cells.Item(2, 7).Value = 0.1022928192822754
$ws.Cells.Item(2, 8).Value = 0.1022928192822754
$ws.Cells.Item(2, 9).Value = 0.1022928192822754
$ws.Cells.Item(2, 10).Value = 0.1022928192822754
$ws.Cells.Item(2, 11).Value = 0.1022928192822754
$ws.Cells.Item(2, 12).Value = 0.1015026103005878
$ws.Cells.Item(2, 13).Value = 0.1015026103005878
$ws.Cells.Item(2, 14).Value = 0.1015026103005878
$ws.Cells.Item(2, 15).Value = 0.1015026103005878
$ws.Cells.Item(2, 16).Value = 0.1015026103005878
$ws.Cells.Item(2, 17).Value = 0.1015026103005878
$ws.Cells.Item(2, 18).Value = 0.1012191830209305
$ws.Cells.Item(2, 19).Value = 0.1012191830209305
$ws.Cells.Item(2, 20).Value = 0.09963959863654148
$ws.Cells.Item(2, 21).Value = 0.09963959863654148
$ws.Cells.Item(2, 22).Value = 0.09963959863654148
$ws.Cells.Item(2, 23).Value = 0.09963959863654148
$ws.Cells.Item(2, 24).Value = 0.09947085426399552
$ws.Cells.Item(2, 25).Value = 0.09906143403966021
$ws.Cells.Item(3, 3).Value = 0.6250035762786865
$ws.Cells.Item(3, 5).Value = 4313.820966464807
$ws.Cells.Item(3, 6).Value = 0.1008992431222507
$ws.Cells.Item(3, 7).Value = 0.1008992431222507
$ws.Cells.Item(3, 8).Value = 0.1008992431222507
$ws.Cells.Item(3, 9).Value = 0.1008992431222507
$ws.Cells.Item(3, 10).Value = 0.1008992431222507
$ws.Cells.Item(3, 11).Value = 0.1008992431222507
$ws.Cells.Item(3, 12).Value = 0.1008992431222507
$ws.Cells.Item(3, 13).Value = 0.1008992431222507
$ws.Cells.Item(3, 14).Value = 0.1008992431222507
$ws.Cells.Item(3, 15).Value = 0.09771982860074843
$ws.Cells.Item(3, 16).Value = 0.09771982860074843
$ws.Cells.Item(3, 17).Value = 0.0956267303660274
$ws.Cells.Item(3, 18).Value = 0.09540436493550206
$ws.Cells.Item(3, 19).Value = 0.09540436493550206
$ws.Cells.Item(3, 20).Value = 0.09540436493550206
$ws.Cells.Item(3, 21).Value = 0.09491316111620605
$ws.Cells.Item(3, 22).Value = 0.09491316111620605
$ws.Cells.Item(3, 23).Value = 0.09491316111620605
$ws.Cells.Item(3, 24).Value = 0.09475674398566873
$ws.Cells.Item(3, 25).Value = 0.09475674398566873
$ws.Cells.Item(4, 3).Value = 0.6679477691650391
$ws.Cells.Item(4, 5).Value = 4578.302601786594
$ws.Cells.Item(4, 6).Value = 0.1014459645610185
$ws.Cells.Item(4, 7).Value = 0.1014459645610185
$ws.Cells.Item(4, 8).Value = 0.1014459645610185
$ws.Cells.Item(4, 9).Value = 0.1014459645610185
$ws.Cells.Item(4, 10).Value = 0.1014459645610185
$ws.Cells.Item(4, 11).Value = 0.1014459645610185
$ws.Cells.Item(4, 12).Value = 0.1014459645610185
$ws.Cells.Item(4, 13).Value = 0.1014459645610185
$ws.Cells.Item(4, 14).Value = 0.1014459645610185
$ws.Cells.Item(4, 15).Value = 0.1014459645610185
$ws.Cells.Item(4, 16).Value = 0.1014459645610185
$ws.Cells.Item(4, 17).Value = 0.1014459645610185
$ws.Cells.Item(4, 18).Value = 0.1014459645610185
$ws.Cells.Item(4, 19).Value = 0.1014459645610185
$ws.Cells.Item(4, 20).Value = 0.1014459645610185
$ws.Cells.Item(4, 21).Value = 0.09991233141884198
$ws.Cells.Item(4, 22).Value = 0.09991233141884198
$ws.Cells.Item(4, 23).Value = 0.09991233141884198
$ws.Cells.Item(4, 24).Value = 0.09991233141884198
$ws.Cells.Item(4, 25).Value = 0.09991233141884198
$ws.Cells.Item(5, 3).Value = 0.5937635898590088
$ws.Cells.Item(5, 5).Value = 4329.091790177639
$ws.Cells.Item(5, 6).Value = 0.1032153233403684
$ws.Cells.Item(5, 7).Value = 0.1032153233403684
$ws.Cells.Item(5, 8).Value = 0.1032153233403684
$ws.Cells.Item(5, 9).Value = 0.1032153233403684
$ws.Cells.Item(5, 10).Value = 0.1032153233403684
$ws.Cells.Item(5, 11).Value = 0.1032153233403684
$ws.Cells.Item(5, 12).Value = 0.1032153233403684
$ws.Cells.Item(5, 13).Value = 0.1032153233403684
$ws.Cells.Item(5, 14).Value = 0.1032153233403684
$ws.Cells.Item(5, 15).Value = 0.09813973000568534
$ws.Cells.Item(5, 16).Value = 0.09813973000568534
$ws.Cells.Item(5, 17).Value = 0.09813973000568534
$ws.Cells.Item(5, 18).Value = 0.0972801682234703
$ws.Cells.Item(5, 19).Value = 0.09505442086116254
$ws.Cells.Item(5, 20).Value = 0.09505442086116254
$ws.Cells.Item(5, 21).Value = 0.09505442086116254
$ws.Cells.Item(5, 22).Value = 0.09505442086116254
$ws.Cells.Item(5, 23).Value = 0.09505442086116254
$ws.Cells.Item(5, 24).Value = 0.09505442086116254
$ws.Cells.Item(5, 25).Value = 0.09505442086116254
$ws.Cells.Item(6, 3).Value = 0.5781099796295166
$ws.Cells.Item(6, 5).Value = 4406.414270384214
$ws.Cells.Item(6, 6).Value = 0.1040766230708546
$ws.Cells.Item(6, 7).Value = 0.1019601353595904
$ws.Cells.Item(6, 8).Value = 0.0984009888803431
$ws.Cells.Item(6, 9).Value = 0.0984009888803431
$ws.Cells.Item(6, 10).Value = 0.0984009888803431
$ws.Cells.Item(6, 11).Value = 0.0984009888803431
$ws.Cells.Item(6, 12).Value = 0.0984009888803431
$ws.Cells.Item(6, 13).Value = 0.0984009888803431
$ws.Cells.Item(6, 14).Value = 0.0984009888803431
$ws.Cells.Item(6, 15).Value = 0.0984009888803431
$ws.Cells.Item(6, 16).Value = 0.0984009888803431
$ws.Cells.Item(6, 17).Value = 0.09700025053487547
$ws.Cells.Item(6, 18).Value = 0.09700025053487547
$ws.Cells.Item(6, 19).Value = 0.09700025053487547
$ws.Cells.Item(6, 20).Value = 0.09700025053487547
$ws.Cells.Item(6, 21).Value = 0.09700025053487547
$ws.Cells.Item(6, 22).Value = 0.09700025053487547
$ws.Cells.Item(6, 23).Value = 0.0968461706015401
$ws.Cells.Item(6, 24).Value = 0.0968461706015401
$ws.Cells.Item(6, 25).Value = 0.09656168168390278
$ws.Cells.Item(7, 3).Value = 0.578113317489624
$ws.Cells.Item(7, 5).Value = 4523.126998402061
$ws.Cells.Item(7, 6).Value = 0.1028491964906056
$ws.Cells.Item(7, 7).Value = 0.1028491964906056
$ws.Cells.Item(7, 8).Value = 0.1028491964906056
$ws.Cells.Item(7, 9).Value = 0.1028491964906056
$ws.Cells.Item(7, 10).Value = 0.1028491964906056
$ws.Cells.Item(7, 11).Value = 0.09963853581770145
$ws.Cells.Item(7, 12).Value = 0.09963853581770145
$ws.Cells.Item(7, 13).Value = 0.09961141509239492
$ws.Cells.Item(7, 14).Value = 0.09961141509239492
$ws.Cells.Item(7, 15).Value = 0.09961141509239492
$ws.Cells.Item(7, 16).Value = 0.09961141509239492
$ws.Cells.Item(7, 17).Value = 0.09961141509239492
$ws.Cells.Item(7, 18).Value = 0.09961141509239492
$ws.Cells.Item(7, 19).Value = 0.09961141509239492
$ws.Cells.Item(7, 20).Value = 0.09961141509239492
$ws.Cells.Item(7, 21).Value = 0.09961141509239492
$ws.Cells.Item(7, 22).Value = 0.09961141509239492
$ws.Cells.Item(7, 23).Value = 0.09961141509239492
$ws.Cells.Item(7, 24).Value = 0.09915056481122853
$ws.Cells.Item(7, 25).Value = 0.09883678359458208
$ws.Cells.Item(8, 3).Value = 0.5781402587890625
$ws.Cells.Item(8, 5).Value = 4471.505068528917
$ws.Cells.Item(8, 6).Value = 0.1021886739605689
$ws.Cells.Item(8, 7).Value = 0.1021886739605689
$ws.Cells.Item(8, 8).Value = 0.1021886739605689
$ws.Cells.Item(8, 9).Value = 0.1021886739605689
$ws.Cells.Item(8, 10).Value = 0.1021886739605689
$ws.Cells.Item(8, 11).Value = 0.1021886739605689
$ws.Cells.Item(8, 12).Value = 0.1021886739605689
$ws.Cells.Item(8, 13).Value = 0.1021886739605689
$ws.Cells.Item(8, 14).Value = 0.1021886739605689
$ws.Cells.Item(8, 15).Value = 0.1008359786319744
$ws.Cells.Item(8, 16).Value = 0.100356936082069
$ws.Cells.Item(8, 17).Value = 0.100356936082069
$ws.Cells.Item(8, 18).Value = 0.1001647421061883
$ws.Cells.Item(8, 19).Value = 0.09894689626702891
$ws.Cells.Item(8, 20).Value = 0.09802711358848874
$ws.Cells.Item(8, 21).Value = 0.0979719219681076
$ws.Cells.Item(8, 22).Value = 0.0979719219681076
$ws.Cells.Item(8, 23).Value = 0.0978544247693687
$ws.Cells.Item(8, 24).Value = 0.0978544247693687
$ws.Cells.Item(8, 25).Value = 0.09783050815845841
$ws.Cells.Item(9, 3).Value = 0.5781095027923584
$ws.Cells.Item(9, 5).Value = 4303.93937970039
$ws.Cells.Item(9, 6).Value = 0.1000250120859812
$ws.Cells.Item(9, 7).Value = 0.1000250120859812
$ws.Cells.Item(9, 8).Value = 0.1000250120859812
$ws.Cells.Item(9, 9).Value = 0.1000250120859812
$ws.Cells.Item(9, 10).Value = 0.1000250120859812
$ws.Cells.Item(9, 11).Value = 0.1000250120859812
$ws.Cells.Item(9, 12).Value = 0.1000250120859812
$ws.Cells.Item(9, 13).Value = 0.09473712315025891
$ws.Cells.Item(9, 14).Value = 0.09473712315025891
$ws.Cells.Item(9, 15).Value = 0.09473712315025891
$ws.Cells.Item(9, 16).Value = 0.09473712315025891
$ws.Cells.Item(9, 17).Value = 0.09473712315025891
$ws.Cells.Item(9, 18).Value = 0.09473712315025891
$ws.Cells.Item(9, 19).Value = 0.09473712315025891
$ws.Cells.Item(9, 20).Value = 0.09473712315025891
$ws.Cells.Item(9, 21).Value = 0.09473712315025891
$ws.Cells.Item(9, 22).Value = 0.0945641204619959
$ws.Cells.Item(9, 23).Value = 0.0945641204619959
$ws.Cells.Item(9, 24).Value = 0.0945641204619959
$ws.Cells.Item(9, 25).Value = 0.0945641204619959
$ws.Cells.Item(10, 3).Value = 0.5937612056732178
$ws.Cells.Item(10, 5).Value = 4366.907621134525
$ws.Cells.Item(10, 6).Value = 0.1045858913633284
$ws.Cells.Item(10, 7).Value = 0.1045858913633284
$ws.Cells.Item(10, 8).Value = 0.1045858913633284
$ws.Cells.Item(10, 9).Value = 0.1045858913633284
$ws.Cells.Item(10, 10).Value = 0.1045858913633284
$ws.Cells.Item(10, 11).Value = 0.1045858913633284
$ws.Cells.Item(10, 12).Value = 0.1045858913633284
$ws.Cells.Item(10, 13).Value = 0.1001760729486956
$ws.Cells.Item(10, 14).Value = 0.09728603818758549
$ws.Cells.Item(10, 15).Value = 0.09728603818758549
$ws.Cells.Item(10, 16).Value = 0.09728603818758549
$ws.Cells.Item(10, 17).Value = 0.09728603818758549
$ws.Cells.Item(10, 18).Value = 0.09728603818758549
$ws.Cells.Item(10, 19).Value = 0.09728603818758549
$ws.Cells.Item(10, 20).Value = 0.09627327579236818
$ws.Cells.Item(10, 21).Value = 0.09627327579236818
$ws.Cells.Item(10, 22).Value = 0.09627327579236818
$ws.Cells.Item(10, 23).Value = 0.09627327579236818
$ws.Cells.Item(10, 24).Value = 0.09610536699934989
$ws.Cells.Item(10, 25).Value = 0.0957915715620765
$ws.Cells.Item(11, 3).Value = 0.5625028610229492
$ws.Cells.Item(11, 5).Value = 4629.566991686879
$ws.Cells.Item(11, 6).Value = 0.1014739318149489
$ws.Cells.Item(11, 7).Value = 0.1014739318149489
$ws.Cells.Item(11, 8).Value = 0.1014739318149489
$ws.Cells.Item(11, 9).Value = 0.1014739318149489
$ws.Cells.Item(11, 10).Value = 0.1014739318149489
$ws.Cells.Item(11, 11).Value = 0.1014739318149489
$ws.Cells.Item(11, 12).Value = 0.1014739318149489
$ws.Cells.Item(11, 13).Value = 0.1014739318149489
$ws.Cells.Item(11, 14).Value = 0.1014739318149489
$ws.Cells.Item(11, 15).Value = 0.1014739318149489
$ws.Cells.Item(11, 16).Value = 0.1014739318149489
$ws.Cells.Item(11, 17).Value = 0.1014739318149489
$ws.Cells.Item(11, 18).Value = 0.1014739318149489
$ws.Cells.Item(11, 19).Value = 0.1014739318149489
$ws.Cells.Item(11, 20).Value = 0.1014739318149489
$ws.Cells.Item(11, 21).Value = 0.1014739318149489
$ws.Cells.Item(11, 22).Value = 0.1014739318149489
$ws.Cells.Item(11, 23).Value = 0.1014739318149489
$ws.Cells.Item(11, 24).Value = 0.1014739318149489
$ws.Cells.Item(11, 25).Value = 0.1009116372648514
